$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row number, new Price (D), new Volume(1h) (E). $null means "unchanged".
$rows = @(
    @{ Row = 2; D = "30.145.33"; E = "  -0.61%  " },
    @{ Row = 3; D = "1.854.96"; E = "  -0.78%  " },
    @{ Row = 4; D = "0.9990"; E = "  -0.14%  " },
    @{ Row = 5; D = "235.38"; E = "  -0.11%  " },
    @{ Row = 6; D = "0.9989"; E = "  -0.14%  " },
    @{ Row = 7; D = $null; E = "  +0.19%  " },
    @{ Row = 8; D = "0.2882"; E = "  +1.31%  " },
    @{ Row = 9; D = $null; E = "  +0.14%  " },
    @{ Row = 10; D = "21.76"; E = "  +1.36%  " },
    @{ Row = 11; D = "0.07965"; E = "  +1.19%  " },
    @{ Row = 12; D = "97.38"; E = "  -0.36%  " },
    @{ Row = 13; D = "1.854.04"; E = "  -0.86%  " },
    @{ Row = 14; D = "5.099"; E = "  +0.15%  " },
    @{ Row = 15; D = "0.6766"; E = "  +0.26%  " },
    @{ Row = 16; D = "267.59"; E = "  -3.24%  " },
    @{ Row = 17; D = "30.129.47"; E = "  -0.62%  " },
    @{ Row = 18; D = "13.61"; E = "  +7.03%  " },
    @{ Row = 19; D = "0.000007638"; E = "  +4.45%  " },
    @{ Row = 20; D = "0.9989"; E = "  -0.14%  " },
    @{ Row = 21; D = "2.096.73"; E = "  -0.55%  " },
    @{ Row = 22; D = "0.9994"; E = "  -0.12%  " },
    @{ Row = 23; D = "5.199"; E = "  -4.92%  " },
    @{ Row = 24; D = "6.141"; E = "  -0.04%  " },
    @{ Row = 25; D = "166.91"; E = "  +0.92%  " },
    @{ Row = 26; D = "9.169"; E = "  +0.47%  " },
    @{ Row = 28; D = "1.937"; E = "  +0.43%  " },
    @{ Row = 29; D = "1.378"; E = "  -0.19%  " },
    @{ Row = 30; D = "0.09875"; E = "  +2.62%  " },
    @{ Row = 31; D = $null; E = "  -0.76%  " },
    @{ Row = 32; D = "4.303"; E = "  -1.68%  " },
    @{ Row = 33; D = "4.021"; E = "  -1.71%  " },
    @{ Row = 34; D = "0.04699"; E = "  -0.13%  " },
    @{ Row = 35; D = $null; E = "  -0.90%  " },
    @{ Row = 36; D = "0.6974"; E = "  -1.18%  " },
    @{ Row = 37; D = "2.703"; E = "  -0.58%  " },
    @{ Row = 38; D = $null; E = "  +0.48%  " },
    @{ Row = 39; D = "2.604"; E = "  +3.03%  " },
    @{ Row = 40; D = "6.315"; E = "  -0.08%  " },
    @{ Row = 41; D = "73.23"; E = "  -0.96%  " },
    @{ Row = 42; D = "1.929"; E = $null },
    @{ Row = 43; D = "0.9982"; E = "  -0.19%  " },
    @{ Row = 44; D = "0.8383"; E = "  -1.31%  " },
    @{ Row = 45; D = "103.30"; E = "  -0.46%  " },
    @{ Row = 46; D = "0.4134"; E = "  -1.25%  " },
    @{ Row = 47; D = "9.138"; E = "  -0.83%  " },
    @{ Row = 48; D = "7.020"; E = "  -2.14%  " },
    @{ Row = 49; D = "932.31"; E = "  -0.48%  " },
    @{ Row = 50; D = "33.90"; E = "  -0.73%  " },
    @{ Row = 51; D = "0.05652"; E = $null }
)

# Cells whose Price column must stay text (trailing zero would be lost as a Number).
$textRows = @(4, 45, 48, 50)

foreach ($item in $rows) {
    if ($null -ne $item.D) {
        $cell = $ws.Cells.Item($item.Row, 4)
        if ($textRows -contains $item.Row) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $item.D
    }
    if ($null -ne $item.E) {
        $ws.Cells.Item($item.Row, 5).Value = $item.E
    }
}